# Usersite Menu - Document Register: add object-repository row for the new
# "Document Register" navigation entry on the Objects_Navigation sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects_Navigation")

# Populate row 7 (directly under the "User Site Menu Navigation" block,
# following the same pattern as the existing rows 3-6).
# Write E7 before B7 so the new shared-string entries land in the same
# order as the target workbook ("Document Register" then
# "Usersite Menu - Document Register").
$ws.Range("E7").Value = "Document Register"
$ws.Range("B7").Value = "Usersite Menu - Document Register"
$ws.Range("C7").Value = "xpath"
$ws.Range("D7").Value = "link"

# B7 keeps the wrapped-text styling used by the other long menu labels
# (e.g. B6 "Usersite Menu - Document & File Storage").
$ws.Range("B7").WrapText = $true

# Move the active selection to the "User Site Sub Menu Transmittals
# Navigation" section header row, matching the saved selection state.
[void]$ws.Range("A10:G10").Select()
